$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    at the top of the data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
# Clone the formatting of the (now shifted) first data row into the
# newly inserted blank row so style (bold index column, borders, etc.)
# matches the rest of the table.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 17
$summary.Range("D2").Value = 1.12

# Renumber the index column for the rows that were pushed down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (so it inherits identical sheet-level formatting / column styles)
#    and inserting it immediately before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$new = $wb.Worksheets.Item("2022-Q3 (2)")
$new.Name = "2022-Q4"

# The template sheet only has 5 data rows (rows 2-6); we need 17 data
# rows (rows 2-18). Extend the sheet by cloning the last data row's
# formatting down through row 18.
$new.Range("A6:H6").Copy()
$new.Range("A7:H18").PasteSpecial(-4122)

$new.Range("A2").Value = 0
$new.Range("B2").Value = "'012284"
$new.Range("C2").Value = "光大保德信健康优加混合"
$new.Range("D2").Value = "'13.09"
$new.Range("E2").Value = "'87.41"
$new.Range("F2").Value = "'4.02"
$new.Range("G2").Value = "'0.5262"
$new.Range("H2").Value = 5
$new.Range("A3").Value = 1
$new.Range("B3").Value = "'519170"
$new.Range("C3").Value = "浦银安盛增长动力灵活配置混合A"
$new.Range("D3").Value = "'6.85"
$new.Range("E3").Value = "'85.43"
$new.Range("F3").Value = "'2.69"
$new.Range("G3").Value = "'0.1843"
$new.Range("H3").Value = 4
$new.Range("A4").Value = 2
$new.Range("B4").Value = "'360005"
$new.Range("C4").Value = "光大保德信红利混合"
$new.Range("D4").Value = "'3.84"
$new.Range("E4").Value = "'79.76"
$new.Range("F4").Value = "'2.66"
$new.Range("G4").Value = "'0.1021"
$new.Range("H4").Value = 7
$new.Range("A5").Value = 3
$new.Range("B5").Value = "'009898"
$new.Range("C5").Value = "民生加银医药健康股票A"
$new.Range("D5").Value = "'2.62"
$new.Range("E5").Value = "'89.21"
$new.Range("F5").Value = "'3.79"
$new.Range("G5").Value = "'0.0993"
$new.Range("H5").Value = 6
$new.Range("A6").Value = 4
$new.Range("B6").Value = "'519113"
$new.Range("C6").Value = "浦银安盛精致生活混合"
$new.Range("D6").Value = "'1.81"
$new.Range("E6").Value = "'90.69"
$new.Range("F6").Value = "'2.43"
$new.Range("G6").Value = "'0.0440"
$new.Range("H6").Value = 9
$new.Range("A7").Value = 5
$new.Range("B7").Value = "'400007"
$new.Range("C7").Value = "东方策略成长混合"
$new.Range("D7").Value = "'1.40"
$new.Range("E7").Value = "'88.54"
$new.Range("F7").Value = "'2.99"
$new.Range("G7").Value = "'0.0419"
$new.Range("H7").Value = 8
$new.Range("A8").Value = 6
$new.Range("B8").Value = "'010690"
$new.Range("C8").Value = "万家互联互通核心资产量化策略混合A"
$new.Range("D8").Value = "'0.53"
$new.Range("E8").Value = "'88.45"
$new.Range("F8").Value = "'6.79"
$new.Range("G8").Value = "'0.0360"
$new.Range("H8").Value = 1
$new.Range("A9").Value = 7
$new.Range("B9").Value = "'014220"
$new.Range("C9").Value = "恒越医疗健康精选混合A"
$new.Range("D9").Value = "'0.73"
$new.Range("E9").Value = "'90.68"
$new.Range("F9").Value = "'3.30"
$new.Range("G9").Value = "'0.0241"
$new.Range("H9").Value = 8
$new.Range("A10").Value = 8
$new.Range("B10").Value = "'530016"
$new.Range("C10").Value = "建信恒稳价值混合"
$new.Range("D10").Value = "'0.54"
$new.Range("E10").Value = "'52.32"
$new.Range("F10").Value = "'3.23"
$new.Range("G10").Value = "'0.0174"
$new.Range("H10").Value = 2
$new.Range("A11").Value = 9
$new.Range("B11").Value = "'006072"
$new.Range("C11").Value = "民生加银创新成长混合A"
$new.Range("D11").Value = "'0.40"
$new.Range("E11").Value = "'91.73"
$new.Range("F11").Value = "'3.16"
$new.Range("G11").Value = "'0.0126"
$new.Range("H11").Value = 8
$new.Range("A12").Value = 10
$new.Range("B12").Value = "'014221"
$new.Range("C12").Value = "恒越医疗健康精选混合C"
$new.Range("D12").Value = "'0.31"
$new.Range("E12").Value = "'90.68"
$new.Range("F12").Value = "'3.30"
$new.Range("G12").Value = "'0.0102"
$new.Range("H12").Value = 8
$new.Range("A13").Value = 11
$new.Range("B13").Value = "'010691"
$new.Range("C13").Value = "万家互联互通核心资产量化策略混合C"
$new.Range("D13").Value = "'0.15"
$new.Range("E13").Value = "'88.45"
$new.Range("F13").Value = "'6.79"
$new.Range("G13").Value = "'0.0102"
$new.Range("H13").Value = 1
$new.Range("A14").Value = 12
$new.Range("B14").Value = "'167703"
$new.Range("C14").Value = "德邦量化优选股票（LOF）C"
$new.Range("D14").Value = "'0.56"
$new.Range("E14").Value = "'88.52"
$new.Range("F14").Value = "'0.97"
$new.Range("G14").Value = "'0.0054"
$new.Range("H14").Value = 8
$new.Range("A15").Value = 13
$new.Range("B15").Value = "'167702"
$new.Range("C15").Value = "德邦量化优选股票（LOF）A"
$new.Range("D15").Value = "'0.32"
$new.Range("E15").Value = "'88.52"
$new.Range("F15").Value = "'0.97"
$new.Range("G15").Value = "'0.0031"
$new.Range("H15").Value = 8
$new.Range("A16").Value = 14
$new.Range("B16").Value = "'014003"
$new.Range("C16").Value = "浦银安盛增长动力灵活配置混合C"
$new.Range("D16").Value = "'0.03"
$new.Range("E16").Value = "'85.43"
$new.Range("F16").Value = "'2.69"
$new.Range("G16").Value = "'0.0008"
$new.Range("H16").Value = 4
$new.Range("A17").Value = 15
$new.Range("B17").Value = "'014758"
$new.Range("C17").Value = "民生加银医药健康股票C"
$new.Range("D17").Value = "'0.01"
$new.Range("E17").Value = "'89.21"
$new.Range("F17").Value = "'3.79"
$new.Range("G17").Value = "'0.0004"
$new.Range("H17").Value = 6
$new.Range("A18").Value = 16
$new.Range("B18").Value = "'014929"
$new.Range("C18").Value = "民生加银创新成长混合C"
$new.Range("D18").Value = "'0.01"
$new.Range("E18").Value = "'91.73"
$new.Range("F18").Value = "'3.16"
$new.Range("G18").Value = "'0.0003"
$new.Range("H18").Value = 8

# Normalise style on the text-typed numeric-looking columns (B, D, E,
# F, G) back to the sheet's default (no bold/border/quote-prefix),
# matching the look of the other data rows; column A keeps the
# index-column style, column H keeps the plain-numeric default style.
$new.Range("B2:B18").Style = "Normal"
$new.Range("D2:G18").Style = "Normal"
